# Swap the contents of rows 13 and 16 (columns A-H) on the active sheet.
# A plain scratch row (far outside the used range) is used as a temporary
# holding area so that Range.Copy can be used for the swap; Copy preserves
# the underlying cell type (shared-string vs number) exactly, which a plain
# Value assignment does not reliably do for numeric-looking text strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = "A100:H100"

# Stash row 13 into the scratch row
$ws.Range("A13:H13").Copy($ws.Range($scratch))

# Move row 16 into row 13
$ws.Range("A16:H16").Copy($ws.Range("A13:H13"))

# Move the stashed original row 13 into row 16
$ws.Range($scratch).Copy($ws.Range("A16:H16"))

# Clear the scratch row (contents and formatting)
$ws.Range($scratch).Clear()

# Row 16 (now holding former row 13's data) had a non-empty D column ("Matriz"),
# which is correct as-is. Row 13 (now holding former row 16's data) must have an
# empty D cell, same as the source row 16 originally had.
$ws.Range("D13").Value = ""
